$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) from A1 onto the three new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add header labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every data row (2 through 48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 100   # column AD
    $ws.Cells.Item($r, 31).Value = 62    # column AE
    $ws.Cells.Item($r, 32).Value = 0     # column AF
}
